$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I6").Value = 12.67039049919483
$ws.Range("N6").Value = 1.847705673092716
$ws.Range("O6").Value = 2.020749367497032
$ws.Range("I7").Value = 12.67039049919483
$ws.Range("I10").Value = 16.86342592592595
$ws.Range("N10").Value = 1.92665172779809
$ws.Range("O10").Value = 2.116885095206829
$ws.Range("I11").Value = 16.86342592592595
$ws.Range("I18").Value = 12.67039049919483
$ws.Range("N18").Value = 1.847705673092716
$ws.Range("O18").Value = 2.020749367497032
$ws.Range("I19").Value = 12.67039049919483
$ws.Range("I29").Value = -1.819444444444444
$ws.Range("N29").Value = 1.618523362263702
$ws.Range("O29").Value = 1.746638928617865
$ws.Range("I30").Value = -1.819444444444444
$ws.Range("I31").Value = 13.62268518518517
$ws.Range("N31").Value = 1.865062221714807
$ws.Range("O31").Value = 2.041808739708676
$ws.Range("I32").Value = 13.62268518518517
$ws.Range("I35").Value = 13.0158303464755
$ws.Range("N35").Value = 1.853964204859962
$ws.Range("O35").Value = 2.02833814451736
$ws.Range("I36").Value = 13.0158303464755
$ws.Range("I37").Value = -1.819444444444444
$ws.Range("N37").Value = 1.618523362263702
$ws.Range("O37").Value = 1.746638928617865
$ws.Range("I38").Value = -1.819444444444444
$ws.Range("I41").Value = 13.0158303464755
$ws.Range("N41").Value = 1.853964204859962
$ws.Range("O41").Value = 2.02833814451736
$ws.Range("I42").Value = 13.0158303464755
$ws.Range("I45").Value = 19.36574074074073
$ws.Range("N45").Value = 1.977063465169192
$ws.Range("O45").Value = 2.178742498783586
$ws.Range("I46").Value = 19.36574074074073
$ws.Range("I47").Value = 13.62268518518517
$ws.Range("N47").Value = 1.865062221714807
$ws.Range("O47").Value = 2.041808739708676
$ws.Range("I48").Value = 13.62268518518517
$ws.Range("I51").Value = 19.65277777777778
$ws.Range("N51").Value = 1.983015294974508
$ws.Range("O51").Value = 2.18606997558991
$ws.Range("I52").Value = 19.65277777777778
$ws.Range("I63").Value = 13.0158303464755
$ws.Range("N63").Value = 1.853964204859962
$ws.Range("O63").Value = 2.02833814451736
$ws.Range("I64").Value = 13.0158303464755
$ws.Range("I67").Value = -1.819444444444444
$ws.Range("N67").Value = 1.618523362263702
$ws.Range("O67").Value = 1.746638928617865
$ws.Range("I68").Value = -1.819444444444444
$ws.Range("I71").Value = 5.486111111111112
$ws.Range("N71").Value = 1.726493341788205
$ws.Range("O71").Value = 1.874863921842289
$ws.Range("I72").Value = 5.486111111111112

Write-Host "Done updating cells"
